# Updates cryptos list prices / volume percentages to the refreshed figures,
# and swaps the Mantle / Bittensor rows (34 and 35), matching the upstream
# data refresh performed by the scheduled GitHub Actions job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that forces a value to be written as plain text, so that strings
# which look numeric (e.g. "597.58", "1.00") are not silently reinterpreted
# by Excel as numbers (which would introduce floating point noise / drop
# trailing zeros). The cell style is restored afterwards so no stray
# "Text" number format lingers on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "70.724.63"
Set-TextValue $ws.Range("E2") "  +2.52%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.585.69"
Set-TextValue $ws.Range("E3") "  +2.07%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("E4") "  -0.03%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "597.58"
Set-TextValue $ws.Range("E5") "  +1.25%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "173.45"
Set-TextValue $ws.Range("E6") "  +1.38%  "

# Row 7 - LidoStakedEther
Set-TextValue $ws.Range("D7") "3.579.28"
Set-TextValue $ws.Range("E7") "  +2.20%  "

# Row 8 - XRP
Set-TextValue $ws.Range("D8") "0.618"
Set-TextValue $ws.Range("E8") "  +0.71%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("E10") "  +6.00%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "7.43"
Set-TextValue $ws.Range("E11") "  +7.42%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.590"
Set-TextValue $ws.Range("E12") "  +1.76%  "

# Row 13 - Avalanche
Set-TextValue $ws.Range("D13") "46.81"
Set-TextValue $ws.Range("E13") "  -0.52%  "

# Row 14 - ShibaInu
Set-TextValue $ws.Range("E14") "  +0.93%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "4.166.35"
Set-TextValue $ws.Range("E15") "  +2.11%  "

# Row 16 - Polkadot
Set-TextValue $ws.Range("D16") "8.42"
Set-TextValue $ws.Range("E16") "  -0.32%  "

# Row 17 - BitcoinCash
Set-TextValue $ws.Range("D17") "613.44"
Set-TextValue $ws.Range("E17") "  -0.99%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.589.07"
Set-TextValue $ws.Range("E18") "  +1.85%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "70.826.96"
Set-TextValue $ws.Range("E19") "  +2.48%  "

# Row 20 - TRON
Set-TextValue $ws.Range("E20") "  -0.85%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "17.48"
Set-TextValue $ws.Range("E21") "  +0.55%  "

# Row 22 - Polygon
Set-TextValue $ws.Range("D22") "0.887"
Set-TextValue $ws.Range("E22") "  +0.31%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "9.31"
Set-TextValue $ws.Range("E23") "  -16.32%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "15.94"
Set-TextValue $ws.Range("E24") "  +0.47%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "97.15"
Set-TextValue $ws.Range("E25") "  +0.38%  "

# Row 26 - PancakeSwap
Set-TextValue $ws.Range("E26") "  -1.93%  "

# Row 27 - Dai
Set-TextValue $ws.Range("D27") "1.00"
Set-TextValue $ws.Range("E27") "  +0.07%  "

# Row 28 - ImmutableX
Set-TextValue $ws.Range("D28") "2.65"
Set-TextValue $ws.Range("E28") "  +0.71%  "

# Row 29 - EthereumClassic
Set-TextValue $ws.Range("D29") "33.91"
Set-TextValue $ws.Range("E29") "  +3.85%  "

# Row 30 - RenderToken
Set-TextValue $ws.Range("D30") "9.19"
Set-TextValue $ws.Range("E30") "  -0.55%  "

# Row 31 - Filecoin
Set-TextValue $ws.Range("D31") "8.43"
Set-TextValue $ws.Range("E31") "  -0.73%  "

# Row 32 - Stacks
Set-TextValue $ws.Range("D32") "3.06"
Set-TextValue $ws.Range("E32") "  -2.08%  "

# Row 33 - NEARProtocol
Set-TextValue $ws.Range("D33") "7.20"
Set-TextValue $ws.Range("E33") "  +4.14%  "

# Rows 34 and 35 swap ranking order: Mantle <-> Bittensor
Set-TextValue $ws.Range("B34") "Bittensor"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D34") "646.64"
Set-TextValue $ws.Range("E34") "  +2.63%  "

Set-TextValue $ws.Range("B35") "Mantle"
Set-TextValue $ws.Range("C35") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D35") "1.30"
Set-TextValue $ws.Range("E35") "  -1.61%  "

# Row 36 - dogwifhat
Set-TextValue $ws.Range("D36") "3.71"
Set-TextValue $ws.Range("E36") "  +6.89%  "

# Row 37 - Hedera
Set-TextValue $ws.Range("D37") "0.101"
Set-TextValue $ws.Range("E37") "  -1.12%  "

# Row 38 - Cosmos
Set-TextValue $ws.Range("D38") "10.85"
Set-TextValue $ws.Range("E38") "  +0.83%  "

# Row 39 - VeChain
Set-TextValue $ws.Range("D39") "0.0482"
Set-TextValue $ws.Range("E39") "  +6.24%  "

# Row 40 - OKB
Set-TextValue $ws.Range("D40") "57.30"
Set-TextValue $ws.Range("E40") "  +0.01%  "

# Row 42 - Kaspa
Set-TextValue $ws.Range("E42") "  +5.10%  "

# Row 43 - Maker
Set-TextValue $ws.Range("D43") "3.400.48"

# Row 44 - TheGraph
Set-TextValue $ws.Range("D44") "0.324"
Set-TextValue $ws.Range("E44") "  -0.77%  "

# Row 45 - PEPE
Set-TextValue $ws.Range("D45") "0.0₃0716"
Set-TextValue $ws.Range("E45") "  +3.19%  "

# Row 46 - InjectiveProtocol
Set-TextValue $ws.Range("D46") "32.98"
Set-TextValue $ws.Range("E46") "  +0.51%  "

# Row 47 - ThetaToken
Set-TextValue $ws.Range("D47") "2.96"
Set-TextValue $ws.Range("E47") "  +6.82%  "

# Row 48 - Fetch.AI
Set-TextValue $ws.Range("D48") "2.66"
Set-TextValue $ws.Range("E48") "  +5.15%  "

# Row 49 - Stellar
Set-TextValue $ws.Range("E49") "  +0.75%  "

# Row 50 - Monero
Set-TextValue $ws.Range("D50") "132.78"
Set-TextValue $ws.Range("E50") "  -0.19%  "

# Row 51 - USDe
Set-TextValue $ws.Range("E51") "  -0.08%  "
